# The workbook's sheet tabs need to be re-sorted: the original left-to-right
# order (2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2022-Q1, 总计) becomes the
# reverse order (总计, 2022-Q1, 2021-Q3, 2021-Q2, 2021-Q1, 2020-Q4).
# No cell data changes - only the tab order (and the sheets' position-derived
# numbering) changes.

$wb = $excel.ActiveWorkbook

# Desired final left-to-right tab order.
$order = @("总计", "2022-Q1", "2021-Q3", "2021-Q2", "2021-Q1", "2020-Q4")

for ($i = 0; $i -lt $order.Length; $i++) {
    $name = $order[$i]
    $ws = $wb.Worksheets.Item($name)
    if ($i -eq 0) {
        # Move the first sheet in the desired order to the very front.
        $ws.Move($wb.Worksheets.Item(1))
    } else {
        # Move each subsequent sheet to immediately follow the one already
        # placed before it.
        $prev = $wb.Worksheets.Item($order[$i - 1])
        $ws.Move($null, $prev)
    }
}
